$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.126.64'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').Value = '2.225.30'
$ws.Range('E3').Value = '  -0.43%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '291.67'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '88.21'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.58%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.514'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.474'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '30.57'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0782'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.72%  '
$ws.Range('E12').Value = '  +3.19%  '
$ws.Range('E13').Value = '  +1.99%  '
$ws.Range('D14').Value = '2.569.24'
$ws.Range('E14').Value = '  -0.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.02'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.83%  '
$ws.Range('D16').Value = '2.203.65'
$ws.Range('E16').Value = '  -1.83%  '
$ws.Range('E17').Value = '  +1.07%  '
$ws.Range('D18').Value = '40.062.63'
$ws.Range('E18').Value = '  +0.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.55'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +7.40%  '
$ws.Range('D20').Value = '0.0₃0887'
$ws.Range('E20').Value = '  -0.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.84'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.81'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.67%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.39'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.61%  '
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.47'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.89%  '
$ws.Range('E26').Value = '  -0.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.73'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.99%  '
$ws.Range('E28').Value = '  -0.87%  '
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '156.05'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '31.85'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.50%  '
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.96'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.24%  '
$ws.Range('E34').Value = '  +1.25%  '
$ws.Range('E35').Value = '  -0.90%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.88'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.24%  '
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '15.83'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.85%  '
$ws.Range('E39').Value = '  -0.32%  '
$ws.Range('E40').Value = '  +2.10%  '
$ws.Range('D41').Value = '2.108.42'
$ws.Range('E41').Value = '  +8.09%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.85'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.37%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.14'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.82%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '18.02'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +10.89%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.00'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.35%  '
$ws.Range('E46').Value = '  -0.64%  '
$ws.Range('E47').Value = '  +2.96%  '
$ws.Range('D48').Value = '2.434.90'
$ws.Range('E48').Value = '  -0.75%  '
$ws.Range('B49').Value = 'TrustWalletToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.11'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.77%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '89.18'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.54%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.45'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.94%  '
